$d = $word.ActiveDocument

# Remove the existing (hidden) _GoBack bookmark; Word will re-create it
# at the location of the next edit.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# Replace " Pseudo-Code" with " Pseudo " + "Code" split across two runs
# with a _GoBack bookmark in between (mirrors Word's own behavior when a
# user edits text in place: the hyphen is replaced by a space).
$r = $d.Range(18, 31)
$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> Pseudo </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Code</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r.InsertXML($xml)
